$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.781.61"
$ws.Range("E2").Value = "  -6.43%  "

# Row 3
$ws.Range("D3").Value = "2.897.14"
$ws.Range("E3").Value = "  -4.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.39%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "2.888.22"
$ws.Range("E8").Value = "  -4.97%  "

# Row 9
$ws.Range("E9").Value = "  -1.54%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.15%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.09%  "

# Row 13
$ws.Range("E13").Value = "  -9.07%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.87%  "

# Row 15
$ws.Range("E15").Value = "  -1.09%  "

# Row 16
$ws.Range("D16").Value = "3.372.45"
$ws.Range("E16").Value = "  -4.88%  "

# Row 17
$ws.Range("D17").Value = "2.890.61"
$ws.Range("E17").Value = "  -5.10%  "

# Row 18
$ws.Range("D18").Value = "57.688.87"
$ws.Range("E18").Value = "  -6.59%  "

# Row 19
$ws.Range("E19").Value = "  +1.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "409.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.652"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.82%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.70%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "

# Row 28
$ws.Range("E28").Value = "  -4.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.95%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0953"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.58%  "

# Row 34
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -13.18%  "

# Row 35
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.896"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.97%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.89%  "

# Row 38
$ws.Range("E38").Value = "  +5.52%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0619"
$ws.Range("E39").Value = "  -11.50%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0343"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.78%  "

# Row 41
$ws.Range("E41").Value = "  -4.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "362.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.19%  "

# Row 43
$ws.Range("D43").Value = "2.610.61"
$ws.Range("E43").Value = "  -2.58%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.54%  "

# Row 45
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.61"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.228"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.34%  "

# Row 48
$ws.Range("E48").Value = "  -1.66%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.69%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.17%  "
